$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the team record columns (Wins/Losses/Ties),
# appended right after the existing "Unnamed: 28" column (AC).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold/centered/bordered style used by the rest of row 1's
# header cells (copy formatting only, from AC1, the neighboring header).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team's 2005 record (71 wins, 91 losses, 0 ties) for every
# player row.
$ws.Range("AD2:AD47").Value = 71
$ws.Range("AE2:AE47").Value = 91
$ws.Range("AF2:AF47").Value = 0
